# Append a new data row (row 61) to the Adafruit IO export sheet, mirroring
# the existing rows' layout: Timestamp | Feed Key | Value | Latitude | Longitude | Elevation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 61

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# Column C ("25") looks numeric; force it to stay text like the other rows
# (which store every value - even numeric-looking ones - as text), then
# drop the temporary number-format style so no stray style index is left
# on the cell.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
